$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 (Marking): Right marks 4 -> 5, Wrong marks -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): Right total 108 -> 135, Wrong total -1 -> -1.2
$ws.Range("B12").Value = 135
$ws.Range("C12").Value = -1.2

# Max display text 107/112 -> 133.8/140
$ws.Range("E12").Value = "133.8/140"
